$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Remove the old "Beck & Rose 2016" entry (row 134, no note yet) -
# it will be re-added further down with its note filled in.
$ws.Rows.Item(134).Delete()

# New literature entries appended at the bottom of the list.
$ws.Range("A141").Value = "Finn & Louviere 1992"
$ws.Range("D141").Value = "original best-worst choice paper"

$ws.Range("D142").Value = "best-worst scaling healthcare"
$ws.Range("A142").Value = "Cheung et al 2016"

$ws.Range("D143").Value = "Best" + [char]8211 + "worst scaling: What it can do for health care research and how to do it"
$ws.Range("A143").Value = "Flynn et al 2007"

$ws.Range("D144").Value = "best worst choice in transportation"
$ws.Range("A144").Value = "Beck & Rose 2016"

$ws.Range("A145").Value = "Hausman & McFadden 1984"
$ws.Range("D145").Value = "multinomal logit cite"

$ws.Range("A146").Value = "marley 2008"
$ws.Range("D146").Value = "check on other authors, but set dependent bw choice"

$ws.Range("A147").Value = "marley 2012"
$ws.Range("D147").Value = "check on other authors, but multi attribute bw choice"

$ws.Range("A148").Value = "flynn 2014"
$ws.Range("D148").Value = "check on other authors, but bw choice review in handbook"

# The sorted range (Data > Sort was previously applied to A2:E141) shrinks by
# one row to match the new row count after the deletion above.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B140"))
$ws.Sort.SetRange($ws.Range("A2:E140"))
$ws.Sort.Apply()

# Window position tweak recorded on re-save.
$excel.ActiveWindow.Left = 15080
$excel.ActiveWindow.Top = 500

# Scroll/selection state left where the author was last working.
$ws.Application.ActiveWindow.ScrollRow = 123
$ws.Range("D149").Select()
